# Epi Info Survey.xlsx template update
# - adds two new list sheets ("Sheet2" = foods eaten, "Sheet4" = yes/no/don't know)
# - adds three new survey questions (symptoms time, foods eaten, hospitalization)
# - renames the "symptoms start date" question's title/description
# - adds "Time" to the DataTypes list
# - tweaks the absolute path / window-size bookkeeping Excel stamps into workbook.xml

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the two new sheets in the right spots so the tab order becomes
#    Sheet1, Sheet4, Sheet2, Sheet3, DataTypes (matching sheetId 5/4/3/2 resp.)
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")

# First inserted sheet grabs the lowest unused default name -> "Sheet2"
$sheetFoods = $wb.Worksheets.Add($null, $sheet1)
# Second inserted sheet grabs the next unused default name -> "Sheet4"
$sheetYesNo = $wb.Worksheets.Add($null, $sheetFoods)
# Put "Sheet4" before "Sheet2" so the tab order matches the target workbook
$sheetYesNo.Move($sheetFoods)

$sheet3 = $wb.Worksheets.Item("Sheet3")
$dataTypes = $wb.Worksheets.Item("DataTypes")

# ---------------------------------------------------------------------------
# 2. Main questionnaire sheet (Sheet1)
# ---------------------------------------------------------------------------

# Row 7 - new "symptoms time" question (title typed first)
$sheet1.Range("A7").Value = "What time did the symptoms start?"

# Row 6 ("When did symptoms start?") - rename the title/description to be
# date-specific now that a matching time question is being added below it.
$sheet1.Range("B6").Value = "Symptoms Date"
$sheet1.Range("C6").Value = "Please enter symptoms date…"

# continue row 7
$sheet1.Range("B7").Value = "Symptoms Time"
$sheet1.Range("C7").Value = "Please enter symptoms time…"
$sheet1.Range("D7").Value = "onse_time"
$sheet1.Range("E7").Value = "Time"
$sheet1.Range("F7").Value = $false

# Row 8 - new "foods eaten" question (checkbox list backed by the new Sheet2 tab)
$sheet1.Range("B8").Value = "Foods Eaten"
$sheet1.Range("D8").Value = "eaten_foods"
$sheet1.Range("G8").Value = "Sheet2"

# Populate the new "Sheet2" tab (foods list) - G8 above points here
$sheetFoods.Range("A1").Value = "Fresh celery"
$sheetFoods.Range("A2").Value = "Grapes"
$sheetFoods.Range("A3").Value = "Peaches"
$sheetFoods.Range("A4").Value = "Apple juice"
$sheetFoods.Range("A5").Value = "Orange juice"

# Populate the new "Sheet4" tab (yes/no/don't know list)
$sheetYesNo.Range("A1").Value = "Yes"
$sheetYesNo.Range("A2").Value = "No"

# Row 9 - new "hospitalized" question (options list backed by the new Sheet4 tab)
$sheet1.Range("G9").Value = "Sheet4"

$sheetYesNo.Range("A3").Value = "Don't know "

# finish filling out rows 8 and 9 on Sheet1
$sheet1.Range("A8").Value = "Select eaten foods:"
$sheet1.Range("C8").Value = "Please select eaten foods…"
$sheet1.Range("E8").Value = "Checkbox"
$sheet1.Range("F8").Value = $false

$sheet1.Range("A9").Value = "Was patient hospitalized?"
$sheet1.Range("B9").Value = "Hospitalization"
$sheet1.Range("D9").Value = "hospitalized"
$sheet1.Range("E9").Value = "Options"
$sheet1.Range("F9").Value = $false

# ---------------------------------------------------------------------------
# 3. DataTypes sheet - add "Time" as a recognised question type
# ---------------------------------------------------------------------------
$dataTypes.Range("A8").Value = "Time"

# ---------------------------------------------------------------------------
# 4. Selections / cursor bookkeeping (best effort, cosmetic)
# ---------------------------------------------------------------------------
$sheetYesNo.Range("D5").Select()
$sheetFoods.Range("E9").Select()
$dataTypes.Range("B9").Select()

$sheet1.Activate()
$sheet1.Range("A1").Select()
